# Add ancestry distribution data to the export files
#
# 1) "Scores" sheet gains three new columns between the existing
#    "Score and results match the original publication" column and the
#    "FTP link"/"License/Terms of Use" columns:
#       - Ancestry Distribution (%) - Source of Variant Associations (GWAS)
#       - Ancestry Distribution (%) - Score Development/Training
#       - Ancestry Distribution (%) - PGS Evaluation
#    The first and third get a value of "European:100" for the single data
#    row; the middle (training) column is left blank.
#
# 2) The "Cohort(s)" value on the "Evaluation Sample Sets" sheet changes its
#    delimiter from ", " to "|".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Scores sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Scores")

# The FTP link hyperlink currently lives on O2; it doesn't automatically
# follow the cell when columns are inserted, so remove it first and
# recreate it afterwards at its new location.
$ws.Range("O2").Hyperlinks.Delete()

# Insert three new, blank columns right before the current "FTP link"
# column (O), pushing FTP link -> R and License/Terms of Use -> S.
$ws.Columns("O:Q").Insert()

# New header row (row 1) values
$ws.Range("O1").Value = "Ancestry Distribution (%) - Source of Variant Associations (GWAS)"
$ws.Range("P1").Value = "Ancestry Distribution (%) - Score Development/Training"
$ws.Range("Q1").Value = "Ancestry Distribution (%) - PGS Evaluation"

# New data row (row 2) values - only the GWAS source and PGS Evaluation
# ancestry distributions are populated for this record.
$ws.Range("O2").Value = "European:100"
$ws.Range("Q2").Value = "European:100"

# Recreate the FTP link hyperlink on its new cell (R2) and restore the
# built-in Hyperlink style that was lost when the cell was re-created by
# the column insert/hyperlink re-add.
$ws.Hyperlinks.Add($ws.Range("R2"), "http://ftp.ebi.ac.uk/pub/databases/spot/pgs/scores/PGS2/ScoringFiles/PGS2.txt.gz")
$ws.Range("R2").Style = "Hyperlink"

# ---------------------------------------------------------------------
# Evaluation Sample Sets sheet
# ---------------------------------------------------------------------
$wsEval = $wb.Worksheets.Item("Evaluation Sample Sets")
$wsEval.Range("P2").Value = "ABC|DEF|KLMN"
